# Update crypto price/volume data as scraped on Mon Sep 25 08:43:37 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '26.193.96'
Set-TextCell 'E2' '  -2.02%  '
Set-TextCell 'D3' '1.581.45'
Set-TextCell 'E3' '  -1.32%  '
Set-TextCell 'D5' '209.55'
Set-TextCell 'E5' '  -1.07%  '
Set-TextCell 'D6' '0.496'
Set-TextCell 'E6' '  -3.19%  '
Set-TextCell 'E8' '  -1.59%  '
Set-TextCell 'D9' '0.245'
Set-TextCell 'E9' '  -0.83%  '
Set-TextCell 'D10' '19.51'
Set-TextCell 'E10' '  -1.15%  '
Set-TextCell 'E11' '  +0.04%  '
Set-TextCell 'D12' '1.804.77'
Set-TextCell 'E12' '  -1.24%  '
Set-TextCell 'B13' 'Polkadot'
Set-TextCell 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D13' '4.05'
Set-TextCell 'E13' '  +0.03%  '
Set-TextCell 'B14' 'WrappedEther'
Set-TextCell 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 'D14' '1.558.49'
Set-TextCell 'E14' '  -3.01%  '
Set-TextCell 'E15' '  -1.50%  '
Set-TextCell 'D16' '64.47'
Set-TextCell 'E16' '  -0.95%  '
Set-TextCell 'D17' '26.199.14'
Set-TextCell 'E17' '  -1.87%  '
Set-TextCell 'D18' '0.0₃0734'
Set-TextCell 'E18' '  -0.93%  '
Set-TextCell 'E19' '  +1.10%  '
Set-TextCell 'E20' '  -0.34%  '
Set-TextCell 'D21' '207.05'
Set-TextCell 'E21' '  -1.61%  '
Set-TextCell 'E23' '  -3.44%  '
Set-TextCell 'D24' '8.88'
Set-TextCell 'E24' '  -1.23%  '
Set-TextCell 'D25' '144.50'
Set-TextCell 'E25' '  +0.60%  '
Set-TextCell 'E26' '  -0.33%  '
Set-TextCell 'E27' '  -1.43%  '
Set-TextCell 'E28' '  -1.64%  '
Set-TextCell 'D29' '15.20'
Set-TextCell 'E29' '  -1.21%  '
Set-TextCell 'E30' '  -1.49%  '
Set-TextCell 'E31' '  -0.83%  '
Set-TextCell 'E32' '  -2.16%  '
Set-TextCell 'D33' '2.94'
Set-TextCell 'E33' '  -1.19%  '
Set-TextCell 'D34' '1.276.06'
Set-TextCell 'E34' '  -1.47%  '
Set-TextCell 'E35' '  -0.36%  '
Set-TextCell 'D36' '0.613'
Set-TextCell 'E36' '  +1.69%  '
Set-TextCell 'E37' '  -1.13%  '
Set-TextCell 'E38' '  -2.09%  '
Set-TextCell 'E39' '  -1.97%  '
Set-TextCell 'D40' '1.01'
Set-TextCell 'E40' '  -12.32%  '
Set-TextCell 'D41' '5.54'
Set-TextCell 'E41' '  +2.43%  '
Set-TextCell 'D42' '0.766'
Set-TextCell 'E42' '  -2.40%  '
Set-TextCell 'E43' '  -2.83%  '
Set-TextCell 'D44' '62.25'
Set-TextCell 'E44' '  -1.32%  '
Set-TextCell 'D45' '1.718.24'
Set-TextCell 'E45' '  -1.12%  '
Set-TextCell 'D46' '89.07'
Set-TextCell 'E46' '  -1.68%  '
Set-TextCell 'E47' '  -0.29%  '
Set-TextCell 'E49' '  -2.06%  '
Set-TextCell 'E50' '  -0.06%  '
Set-TextCell 'D51' '7.41'
Set-TextCell 'E51' '  -0.29%  '
